$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2026-02-04 17:15:04"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2026-02-04 17:15:07"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "90%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2026-02-04 17:15:09"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "83%"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "993.4 hPa"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "7.8 MJ/m2"
$ws.Range("O4").NumberFormat = "@"
$ws.Range("O4").Value = "5.7 °C"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2026-02-04 17:15:12"
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = "992.8 hPa"
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = "9.5 MJ/m2"
$ws.Range("O5").NumberFormat = "@"
$ws.Range("O5").Value = "8.5 °C"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2026-02-04 17:15:15"
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "993.9 hPa"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "5.3 MJ/m2"
$ws.Range("O6").NumberFormat = "@"
$ws.Range("O6").Value = "10.9 °C"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2026-02-04 17:15:17"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "5.6 MJ/m2"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2026-02-04 17:15:20"
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "8.5 MJ/m2"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2026-02-04 17:15:22"
$ws.Range("O9").NumberFormat = "@"
$ws.Range("O9").Value = "2.1 °C"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2026-02-04 17:15:24"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "86%"
$ws.Range("O10").NumberFormat = "@"
$ws.Range("O10").Value = "8.2 °C"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2026-02-04 17:15:27"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "83%"
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = "995.3 hPa"
$ws.Range("K11").NumberFormat = "@"
$ws.Range("K11").Value = "10.6 MJ/m2"
$ws.Range("O11").NumberFormat = "@"
$ws.Range("O11").Value = "0.5 °C"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2026-02-04 17:15:29"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2026-02-04 17:15:32"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "81%"
$ws.Range("O13").NumberFormat = "@"
$ws.Range("O13").Value = "7.2 °C"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2026-02-04 17:15:34"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "83%"
$ws.Range("K14").NumberFormat = "@"
$ws.Range("K14").Value = "9.2 MJ/m2"
$ws.Range("O14").NumberFormat = "@"
$ws.Range("O14").Value = "-6.4 °C"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2026-02-04 17:15:37"
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = "993.5 hPa"
$ws.Range("O15").NumberFormat = "@"
$ws.Range("O15").Value = "6.2 °C"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2026-02-04 17:15:39"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "92%"
$ws.Range("O16").NumberFormat = "@"
$ws.Range("O16").Value = "2.8 °C"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2026-02-04 17:15:42"
$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value = "994.7 hPa"
$ws.Range("O17").NumberFormat = "@"
$ws.Range("O17").Value = "3.1 °C"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2026-02-04 17:15:44"
$ws.Range("K18").NumberFormat = "@"
$ws.Range("K18").Value = "6.3 MJ/m2"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2026-02-04 17:15:47"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2026-02-04 17:15:49"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "113 cm"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2026-02-04 17:15:52"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "74%"
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value = "993.0 hPa"
$ws.Range("O21").NumberFormat = "@"
$ws.Range("O21").Value = "6.0 °C"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2026-02-04 17:15:55"
$ws.Range("K22").NumberFormat = "@"
$ws.Range("K22").Value = "8.3 MJ/m2"
$ws.Range("O22").NumberFormat = "@"
$ws.Range("O22").Value = "8.3 °C"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2026-02-04 17:15:57"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "71%"
$ws.Range("J23").NumberFormat = "@"
$ws.Range("J23").Value = "992.4 hPa"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2026-02-04 17:16:00"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "71%"
$ws.Range("J24").NumberFormat = "@"
$ws.Range("J24").Value = "991.5 hPa"
$ws.Range("K24").NumberFormat = "@"
$ws.Range("K24").Value = "9.5 MJ/m2"
$ws.Range("O24").NumberFormat = "@"
$ws.Range("O24").Value = "10.0 °C"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2026-02-04 17:16:02"
$ws.Range("J25").NumberFormat = "@"
$ws.Range("J25").Value = "994.2 hPa"
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value = "10.7 MJ/m2"
$ws.Range("O25").NumberFormat = "@"
$ws.Range("O25").Value = "1.1 °C"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2026-02-04 17:16:04"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2026-02-04 17:16:07"
$ws.Range("J27").NumberFormat = "@"
$ws.Range("J27").Value = "993.1 hPa"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "2026-02-04 17:16:09"
$ws.Range("J28").NumberFormat = "@"
$ws.Range("J28").Value = "994.9 hPa"
$ws.Range("O28").NumberFormat = "@"
$ws.Range("O28").Value = "2.3 °C"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "2026-02-04 17:16:11"
$ws.Range("O29").NumberFormat = "@"
$ws.Range("O29").Value = "7.1 °C"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "2026-02-04 17:16:14"
$ws.Range("O30").NumberFormat = "@"
$ws.Range("O30").Value = "-5.5 °C"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "2026-02-04 17:16:16"
$ws.Range("J31").NumberFormat = "@"
$ws.Range("J31").Value = "994.8 hPa"
$ws.Range("O31").NumberFormat = "@"
$ws.Range("O31").Value = "4.1 °C"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "2026-02-04 17:16:19"
$ws.Range("J32").NumberFormat = "@"
$ws.Range("J32").Value = "993.6 hPa"
$ws.Range("O32").NumberFormat = "@"
$ws.Range("O32").Value = "10.3 °C"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "2026-02-04 17:16:21"
$ws.Range("O33").NumberFormat = "@"
$ws.Range("O33").Value = "9.7 °C"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "2026-02-04 17:16:23"
$ws.Range("K34").NumberFormat = "@"
$ws.Range("K34").Value = "7.1 MJ/m2"
$ws.Range("O34").NumberFormat = "@"
$ws.Range("O34").Value = "3.3 °C"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "2026-02-04 17:16:26"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "2026-02-04 17:16:28"
$ws.Range("O36").NumberFormat = "@"
$ws.Range("O36").Value = "7.0 °C"
